$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.8838797317479079
$ws.Range("A3").Value = 0.5804756250896509
$ws.Range("A4").Value = 0.5689877857268043
$ws.Range("A5").Value = 0.5560926013731577
$ws.Range("A6").Value = 0.6009613394303588
$ws.Range("A7").Value = 0.7116331327055651
$ws.Range("A8").Value = 0.678844926177689
$ws.Range("A9").Value = 0.4874762401984203
$ws.Range("A10").Value = 0.4454385935230377
$ws.Range("A11").Value = 0.3114347509580857
$ws.Range("A12").Value = 0.3349623807193612
$ws.Range("A13").Value = 0.1746741661377086
$ws.Range("A14").Value = 0.2179316971256372
$ws.Range("A15").Value = 0.1859517758672763
$ws.Range("A16").Value = 0.1581025253002655
$ws.Range("A17").Value = 0.1287986758911418
$ws.Range("A18").Value = 0.1072357221515841
$ws.Range("A19").Value = 0.09948603012936319
$ws.Range("A20").Value = 0.09377771078667768
$ws.Range("A21").Value = 0.093224478992346
$ws.Range("A22").Value = 0.08882978148549628
$ws.Range("A23").Value = 0.08858154975821175
$ws.Range("A24").Value = 0.08909639842212708
$ws.Range("A25").Value = 0.08806155458954379
$ws.Range("A26").Value = 0.0856825622416382
$ws.Range("A27").Value = 0.08596083630439189
$ws.Range("A28").Value = 0.08520685196275338
$ws.Range("A29").Value = 0.08515360638446175
$ws.Range("A30").Value = 0.08510105188471401
$ws.Range("A31").Value = 0.0849349134539981
$ws.Range("A32").Value = 0.08492314055169457
$ws.Range("A33").Value = 0.08484652061346608
$ws.Range("A34").Value = 0.0848629246688912
$ws.Range("A35").Value = 0.08483236540724776
$ws.Range("A36").Value = 0.08477981081747249
$ws.Range("A37").Value = 0.08475837594653182
$ws.Range("A38").Value = 0.08475730470400228
$ws.Range("A39").Value = 0.08474557933377681
$ws.Range("A40").Value = 0.08471829387392681
$ws.Range("A41").Value = 0.08474114510300992
$ws.Range("A42").Value = 0.08476343409439994
$ws.Range("A43").Value = 0.08472204855566307
$ws.Range("A44").Value = 0.08466100834992524
$ws.Range("A45").Value = 0.08464032251196726
$ws.Range("A46").Value = 0.08460949939466041
$ws.Range("A47").Value = 0.08458676911291615
$ws.Range("A48").Value = 0.08456284622344548
$ws.Range("A49").Value = 0.08454967045776435
$ws.Range("A50").Value = 0.08453807217629299
$ws.Range("A51").Value = 0.0845354999738257
